$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.3067
$ws.Range("E4").Value = 16.07599999999999
$ws.Range("D7").Value = -7.085800000000001
$ws.Range("C8").Value = -12.73429999999999
$ws.Range("C10").Value = -12.98649999999999
$ws.Range("E11").Value = 16.46449999999999
$ws.Range("C12").Value = -10.66909999999999
$ws.Range("D14").Value = -8.515400000000005
$ws.Range("E14").Value = 16.41909999999999
$ws.Range("D15").Value = -7.780400000000001
$ws.Range("C18").Value = -13.47720000000001
$ws.Range("D18").Value = -9.228899999999987
$ws.Range("E18").Value = 16.2467
$ws.Range("E19").Value = 16.61990000000001
$ws.Range("D20").Value = -7.955499999999994
$ws.Range("E21").Value = 17.02030000000001
$ws.Range("C25").Value = -13.39589999999999
$ws.Range("E27").Value = 16.6218
$ws.Range("D29").Value = -6.994499999999999
$ws.Range("D30").Value = -7.866400000000006
$ws.Range("D31").Value = -7.662600000000001
$ws.Range("E31").Value = 16.74500000000001
$ws.Range("D35").Value = -8.375099999999996
$ws.Range("C37").Value = -13.6497
$ws.Range("E38").Value = 16.4403
$ws.Range("D40").Value = -8.141099999999996
$ws.Range("E42").Value = 16.26869999999999
$ws.Range("D44").Value = -7.595499999999998
$ws.Range("E44").Value = 16.4952
$ws.Range("E47").Value = 16.5521
$ws.Range("D50").Value = -7.834699999999994
$ws.Range("D54").Value = -8.371100000000004
$ws.Range("C55").Value = -13.81339999999999
$ws.Range("E56").Value = 16.2754
$ws.Range("E58").Value = 16.24730000000001
$ws.Range("E65").Value = 17.2982
$ws.Range("C68").Value = -10.5693
$ws.Range("D68").Value = -7.042599999999996
$ws.Range("E73").Value = 17.40710000000001
$ws.Range("D76").Value = -7.6109
$ws.Range("C77").Value = -12.9213
$ws.Range("C78").Value = -13.40320000000001
$ws.Range("C79").Value = -12.2474
$ws.Range("C80").Value = -13.5377
$ws.Range("C81").Value = -12.9777
$ws.Range("C82").Value = -12.43679999999999
$ws.Range("C84").Value = -13.5726
$ws.Range("D87").Value = -7.857499999999999
$ws.Range("D88").Value = -7.181999999999997
$ws.Range("E90").Value = 16.3862
$ws.Range("D92").Value = -6.971399999999996
$ws.Range("E92").Value = 17.62330000000002
$ws.Range("E94").Value = 18.98780000000002
$ws.Range("E95").Value = 18.05510000000002
$ws.Range("D96").Value = -8.169900000000004
$ws.Range("D98").Value = -8.440900000000006
$ws.Range("C101").Value = -12.46799999999999
$ws.Range("D101").Value = -8.284199999999997
$ws.Range("E101").Value = 16.5025
$ws.Range("C102").Value = -13.8409
$ws.Range("D102").Value = -7.752800000000001
